# "Adicionei o Join Tabs" - update the aggregated "copy" column (E) with the
# new example labels used by the Join Tabs feature, replacing the old
# placeholder text ("Copia"/"truta"/"tias") with "Exemple"/"Exemple2"/"Exemple3".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Block 1 (rows 2-5, under "thiago"/"bruno" group) -> Exemple
$ws.Range("E2:E5").Value = "Exemple"

# Block 2 (rows 6-9, under "bruno"/"renan" group) -> Exemple2
$ws.Range("E6:E9").Value = "Exemple2"

# Block 3 (rows 10-17, under "Rita"/"Josias" group) -> Exemple3
$ws.Range("E10:E17").Value = "Exemple3"

# Match the author's final cursor position / view state.
$ws.Range("F13").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
